$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Helper: find a paragraph index (1-based) whose Range.Text starts with the
# given literal prefix. Re-scanned every time since indices shift as we
# merge/replace paragraphs.
# ---------------------------------------------------------------------------
function Find-ParaIndex($prefix) {
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($p.Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Hunk 1: merge "Fix the multi-line exceptions." (ilvl0) and the following
# "Test it propagates in raw representation." (ilvl1) paragraph into a
# single ilvl0 paragraph reading "lg help gets stuck." + line break +
# the _GoBack bookmark.
# ---------------------------------------------------------------------------
$iA = Find-ParaIndex("Fix the multi-line exceptions.")
$iB = $iA + 1
$pA = $d.Paragraphs($iA)
$pB = $d.Paragraphs($iB)
$rng1 = $d.Range($pA.Range.Start, $pB.Range.End)

$xml1 = "<w:p $wns>" +
    "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"10`"/></w:numPr><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr></w:pPr>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t>lg</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t xml:space=`"preserve`"> help</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t>gets stuck.</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:br/></w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
    "</w:p>"

[void]$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Hunk 2: the "Fix the off-line properties." (ilvl0) paragraph, the
# following "Test it propagates in raw representation." (ilvl1) paragraph,
# and the (now-orphan, bookmark-only) paragraph after it collapse into a
# single ilvl0 paragraph reading "User manual that mentions -f" + a
# trailing line break.
# ---------------------------------------------------------------------------
$iC = Find-ParaIndex("Fix the off-line properties.")
$iE = $iC + 2
$pC = $d.Paragraphs($iC)
$pE = $d.Paragraphs($iE)
$rng2 = $d.Range($pC.Range.Start, $pE.Range.End)

$xml2 = "<w:p $wns>" +
    "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"10`"/></w:numPr><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t>User manual that mentions -f</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:br/></w:r>" +
    "</w:p>"

[void]$rng2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Hunk 3: append a trailing line break to the
# "Consolidate findNext() across components." paragraph.
# ---------------------------------------------------------------------------
$iF = Find-ParaIndex("Consolidate ")
$pF = $d.Paragraphs($iF)
$rng3 = $d.Range($pF.Range.Start, $pF.Range.End)

$xml3 = "<w:p $wns>" +
    "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"10`"/></w:numPr><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t xml:space=`"preserve`">Consolidate </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t>findNext</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t xml:space=`"preserve`">() across </w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t>components</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:t>.</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Garamond`" w:hAnsi=`"Garamond`"/></w:rPr><w:br/></w:r>" +
    "</w:p>"

[void]$rng3.InsertXML($xml3)

Write-Output "done"
